$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing threshold values (Min/Max columns) ---
$ws.Range("B2").Value = 5.5
$ws.Range("C2").Value = 10.5
$ws.Range("B3").Value = 5
$ws.Range("C5").Value = 15

# --- Append a new (currently empty) row 6, inheriting row 5's formatting ---
$ws.Range("A5:C5").Copy()
$ws.Range("A6:C6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# --- Resize columns (closest attainable widths to the authored sizes) ---
$ws.Columns("A").ColumnWidth = 20.714285714285715
$ws.Columns("B").ColumnWidth = 4.428571428571429
$ws.Columns("C").ColumnWidth = 4.714285714285714
